$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new rows to make room for "agent" (new row 2) and "developer" (new row 7)
$ws.Rows(2).Insert()
$ws.Rows(7).Insert()

# Copy cell formatting (style) from row 3 into the two newly inserted blank rows (2 and 7)
$ws.Range("A3:Y3").Copy()
$ws.Range("A2:Y2").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:Y3").Copy()
$ws.Range("A7:Y7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Step 2: write every cell value for rows 2-20 to match the target layout
# Row 2: agent
$ws.Range("A2").Value = 'BCIO:050273'
$ws.Range("B2").Value = 'agent'
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = 'An independent continuant that is a human being, group or organisation.'
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 'independent continuant'
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = 'Intervention development'
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = ""
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = ""
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = 'RW'
$ws.Range("W2").Value = 'Proposed'
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = ""

# Row 3: behaviour change intervention development process 
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = 'behaviour change intervention development process '
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 'intervention development process'
$ws.Range("H3").Value = 'process'
$ws.Range("I3").Value = 'Intervention development'
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = ""
$ws.Range("R3").Value = ""
$ws.Range("S3").Value = ""
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = '0'
$ws.Range("U3").Value = ""
$ws.Range("V3").Value = 'RW'
$ws.Range("W3").Value = 'Proposed'
$ws.Range("X3").Value = ""
$ws.Range("Y3").Value = ""

# Row 4: co-production
$ws.Range("A4").Value = 'BCIO:050269'
$ws.Range("B4").Value = 'co-production'
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = 'A development process in which developers involve other stakeholders as development partners.'
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = 'development process'
$ws.Range("H4").Value = 'process'
$ws.Range("I4").Value = 'Intervention development'
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = ""
$ws.Range("Q4").Value = ""
$ws.Range("R4").Value = ""
$ws.Range("S4").Value = ""
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value = '0'
$ws.Range("U4").Value = ""
$ws.Range("V4").Value = 'RW'
$ws.Range("W4").Value = 'Proposed'
$ws.Range("X4").Value = ""
$ws.Range("Y4").Value = ""

# Row 5: collaboration process
$ws.Range("A5").Value = 'BCIO:050270'
$ws.Range("B5").Value = 'collaboration process'
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 'A process'
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = 'development process'
$ws.Range("H5").Value = 'process'
$ws.Range("I5").Value = 'Intervention development'
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = ""
$ws.Range("Q5").Value = ""
$ws.Range("R5").Value = ""
$ws.Range("S5").Value = ""
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value = '0'
$ws.Range("U5").Value = ""
$ws.Range("V5").Value = 'RW'
$ws.Range("W5").Value = 'Proposed'
$ws.Range("X5").Value = ""
$ws.Range("Y5").Value = ""

# Row 6: consultation
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = 'consultation'
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 'development process'
$ws.Range("H6").Value = 'process'
$ws.Range("I6").Value = 'Intervention development'
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = ""
$ws.Range("Q6").Value = ""
$ws.Range("R6").Value = ""
$ws.Range("S6").Value = ""
$ws.Range("T6").NumberFormat = "@"
$ws.Range("T6").Value = '0'
$ws.Range("U6").Value = ""
$ws.Range("V6").Value = 'RW'
$ws.Range("W6").Value = 'Proposed'
$ws.Range("X6").Value = ""
$ws.Range("Y6").Value = ""

# Row 7: developer
$ws.Range("A7").Value = 'BCIO:050274'
$ws.Range("B7").Value = 'developer'
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = 'An agent with a developer role.'
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = 'agent'
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = 'Intervention development'
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = ""
$ws.Range("P7").Value = ""
$ws.Range("Q7").Value = ""
$ws.Range("R7").Value = ""
$ws.Range("S7").Value = ""
$ws.Range("T7").Value = ""
$ws.Range("U7").Value = ""
$ws.Range("V7").Value = 'RW'
$ws.Range("W7").Value = 'Proposed'
$ws.Range("X7").Value = ""
$ws.Range("Y7").Value = ""

# Row 8: developer role
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = 'developer role'
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = 'role'
$ws.Range("H8").Value = 'process'
$ws.Range("I8").Value = 'Intervention development'
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("O8").Value = ""
$ws.Range("P8").Value = ""
$ws.Range("Q8").Value = ""
$ws.Range("R8").Value = ""
$ws.Range("S8").Value = ""
$ws.Range("T8").NumberFormat = "@"
$ws.Range("T8").Value = '0'
$ws.Range("U8").Value = ""
$ws.Range("V8").Value = 'RW'
$ws.Range("W8").Value = 'Proposed'
$ws.Range("X8").Value = ""
$ws.Range("Y8").Value = ""

# Row 9: development partner
$ws.Range("A9").Value = 'BCIO:050275'
$ws.Range("B9").Value = 'development partner'
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = 'A person, group or organisation that has a development partner role.'
$ws.Range("E9").Value = '(''human being'' or ''group'' or ''organisation'') and (''has role'' ''development partner role'')'
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = 'agent'
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = 'Intervention development'
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = 'development partner role'
$ws.Range("N9").Value = ""
$ws.Range("O9").Value = ""
$ws.Range("P9").Value = ""
$ws.Range("Q9").Value = ""
$ws.Range("R9").Value = ""
$ws.Range("S9").Value = ""
$ws.Range("T9").NumberFormat = "@"
$ws.Range("T9").Value = '0'
$ws.Range("U9").Value = ""
$ws.Range("V9").Value = 'RW'
$ws.Range("W9").Value = 'Proposed'
$ws.Range("X9").Value = ""
$ws.Range("Y9").Value = ""

# Row 10: development partner role
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = 'development partner role'
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = 'role'
$ws.Range("H10").Value = 'independent continuant'
$ws.Range("I10").Value = 'Intervention development'
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("O10").Value = ""
$ws.Range("P10").Value = ""
$ws.Range("Q10").Value = ""
$ws.Range("R10").Value = ""
$ws.Range("S10").Value = ""
$ws.Range("T10").NumberFormat = "@"
$ws.Range("T10").Value = '0'
$ws.Range("U10").Value = ""
$ws.Range("V10").Value = 'RW'
$ws.Range("W10").Value = 'Proposed'
$ws.Range("X10").Value = ""
$ws.Range("Y10").Value = ""

# Row 11: development process
$ws.Range("A11").Value = 'BCIO:050271'
$ws.Range("B11").Value = 'development process'
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = 'A process that is creation a product, commodity, service or intervention.'
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = 'planned process'
$ws.Range("H11").Value = 'process'
$ws.Range("I11").Value = 'Intervention development'
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("O11").Value = ""
$ws.Range("P11").Value = ""
$ws.Range("Q11").Value = ""
$ws.Range("R11").Value = ""
$ws.Range("S11").Value = ""
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = '0'
$ws.Range("U11").Value = ""
$ws.Range("V11").Value = 'RW'
$ws.Range("W11").Value = 'Proposed'
$ws.Range("X11").Value = ""
$ws.Range("Y11").Value = ""

# Row 12: engagement
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = 'engagement'
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = 'development process'
$ws.Range("H12").Value = 'process'
$ws.Range("I12").Value = 'Intervention development'
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("P12").Value = ""
$ws.Range("Q12").Value = ""
$ws.Range("R12").Value = ""
$ws.Range("S12").Value = ""
$ws.Range("T12").NumberFormat = "@"
$ws.Range("T12").Value = '0'
$ws.Range("U12").Value = ""
$ws.Range("V12").Value = 'RW'
$ws.Range("W12").Value = 'Proposed'
$ws.Range("X12").Value = ""
$ws.Range("Y12").Value = ""

# Row 13: intervention development process
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = 'intervention development process'
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = 'development process'
$ws.Range("H13").Value = 'process'
$ws.Range("I13").Value = 'Intervention development'
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("O13").Value = ""
$ws.Range("P13").Value = ""
$ws.Range("Q13").Value = ""
$ws.Range("R13").Value = ""
$ws.Range("S13").Value = ""
$ws.Range("T13").NumberFormat = "@"
$ws.Range("T13").Value = '0'
$ws.Range("U13").Value = ""
$ws.Range("V13").Value = 'RW'
$ws.Range("W13").Value = 'Proposed'
$ws.Range("X13").Value = ""
$ws.Range("Y13").Value = ""

# Row 14: patient and public involvement
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = 'patient and public involvement'
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = 'development process'
$ws.Range("H14").Value = 'process'
$ws.Range("I14").Value = 'Intervention development'
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("O14").Value = ""
$ws.Range("P14").Value = ""
$ws.Range("Q14").Value = ""
$ws.Range("R14").Value = ""
$ws.Range("S14").Value = ""
$ws.Range("T14").NumberFormat = "@"
$ws.Range("T14").Value = '0'
$ws.Range("U14").Value = ""
$ws.Range("V14").Value = 'RW'
$ws.Range("W14").Value = 'Proposed'
$ws.Range("X14").Value = ""
$ws.Range("Y14").Value = ""

# Row 15: patient and public involvement and engagement
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = 'patient and public involvement and engagement'
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = 'development process'
$ws.Range("H15").Value = 'process'
$ws.Range("I15").Value = 'Intervention development'
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("O15").Value = ""
$ws.Range("P15").Value = ""
$ws.Range("Q15").Value = ""
$ws.Range("R15").Value = ""
$ws.Range("S15").Value = ""
$ws.Range("T15").NumberFormat = "@"
$ws.Range("T15").Value = '0'
$ws.Range("U15").Value = ""
$ws.Range("V15").Value = 'RW'
$ws.Range("W15").Value = 'Proposed'
$ws.Range("X15").Value = ""
$ws.Range("Y15").Value = ""

# Row 16: product development process
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = 'product development process'
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = 'development process'
$ws.Range("H16").Value = 'process'
$ws.Range("I16").Value = 'Intervention development'
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("O16").Value = ""
$ws.Range("P16").Value = ""
$ws.Range("Q16").Value = ""
$ws.Range("R16").Value = ""
$ws.Range("S16").Value = ""
$ws.Range("T16").NumberFormat = "@"
$ws.Range("T16").Value = '0'
$ws.Range("U16").Value = ""
$ws.Range("V16").Value = 'RW'
$ws.Range("W16").Value = 'Proposed'
$ws.Range("X16").Value = ""
$ws.Range("Y16").Value = ""

# Row 17: project development process
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = 'project development process'
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = 'development process'
$ws.Range("H17").Value = 'process'
$ws.Range("I17").Value = 'Intervention development'
$ws.Range("J17").Value = ""
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("O17").Value = ""
$ws.Range("P17").Value = ""
$ws.Range("Q17").Value = ""
$ws.Range("R17").Value = ""
$ws.Range("S17").Value = ""
$ws.Range("T17").NumberFormat = "@"
$ws.Range("T17").Value = '0'
$ws.Range("U17").Value = ""
$ws.Range("V17").Value = 'RW'
$ws.Range("W17").Value = 'Proposed'
$ws.Range("X17").Value = ""
$ws.Range("Y17").Value = ""

# Row 18: service development process
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = 'service development process'
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = 'development process'
$ws.Range("H18").Value = 'process'
$ws.Range("I18").Value = 'Intervention development'
$ws.Range("J18").Value = ""
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = ""
$ws.Range("O18").Value = ""
$ws.Range("P18").Value = ""
$ws.Range("Q18").Value = ""
$ws.Range("R18").Value = ""
$ws.Range("S18").Value = ""
$ws.Range("T18").NumberFormat = "@"
$ws.Range("T18").Value = '0'
$ws.Range("U18").Value = ""
$ws.Range("V18").Value = 'RW'
$ws.Range("W18").Value = 'Proposed'
$ws.Range("X18").Value = ""
$ws.Range("Y18").Value = ""

# Row 19: stakeholder
$ws.Range("A19").Value = 'BCIO:050276'
$ws.Range("B19").Value = 'stakeholder'
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = 'A person, group or organisation that has a stakeholder role.'
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = 'agent'
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = 'Intervention development'
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = 'stakeholder role'
$ws.Range("N19").Value = ""
$ws.Range("O19").Value = ""
$ws.Range("P19").Value = ""
$ws.Range("Q19").Value = ""
$ws.Range("R19").Value = ""
$ws.Range("S19").Value = ""
$ws.Range("T19").NumberFormat = "@"
$ws.Range("T19").Value = '0'
$ws.Range("U19").Value = ""
$ws.Range("V19").Value = 'RW'
$ws.Range("W19").Value = 'Proposed'
$ws.Range("X19").Value = ""
$ws.Range("Y19").Value = ""

# Row 20: stakeholder role
$ws.Range("A20").Value = 'BCIO:050272'
$ws.Range("B20").Value = 'stakeholder role'
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = 'A role that is being involved with or affected by a project, service, intervention, commodity, product or enterprise.'
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = 'role'
$ws.Range("H20").Value = 'independent continuant'
$ws.Range("I20").Value = 'Intervention development'
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""
$ws.Range("O20").Value = ""
$ws.Range("P20").Value = ""
$ws.Range("Q20").Value = ""
$ws.Range("R20").Value = ""
$ws.Range("S20").Value = ""
$ws.Range("T20").NumberFormat = "@"
$ws.Range("T20").Value = '0'
$ws.Range("U20").Value = ""
$ws.Range("V20").Value = 'RW'
$ws.Range("W20").Value = 'Proposed'
$ws.Range("X20").Value = ""
$ws.Range("Y20").Value = ""
